# Apply the "simulator full-month coverage, persist logs, fix employees" edit.
#
# 1. Fix employee/client names in the "Weekly Timesheet" sheet (B2:B6) that
#    were scrambled by the (now fixed) simulator, and correct the employee ID.
# 2. Populate the previously-zeroed Rate/Total columns on both sheets now
#    that the simulator persists its computed logs.

$wb = $excel.ActiveWorkbook

$tsSheet = $wb.Worksheets.Item("Weekly Timesheet")
$jsSheet = $wb.Worksheets.Item("Jason Schema")

# --- Fix scrambled client names (shared between both sheets - the Weekly
#     Timesheet's Client column and the Jason Schema's Client column both
#     point at the same underlying string per day) ---
$tsSheet.Range("B2").Value = "Tubergen"
$tsSheet.Range("B3").Value = "Hewett"
$tsSheet.Range("B4").Value = "Durfee"
$tsSheet.Range("B5").Value = "Markfield"
$tsSheet.Range("B6").Value = "Corr"

$jsSheet.Range("D2").Value = "Tubergen"
$jsSheet.Range("D3").Value = "Hewett"
$jsSheet.Range("D4").Value = "Durfee"
$jsSheet.Range("D5").Value = "Markfield"
$jsSheet.Range("D6").Value = "Corr"

# --- Weekly Timesheet: populate Rate (E) / Total (F) for the daily rows ---
for ($r = 2; $r -le 6; $r++) {
    $tsSheet.Cells.Item($r, 5).Value = 150
    $tsSheet.Cells.Item($r, 6).Value = 1200
}

# --- Weekly Timesheet: subtotal rows now reflect the populated totals ---
$tsSheet.Range("F8").Value = 6000
$tsSheet.Range("F12").Value = 6000
$tsSheet.Range("F13").Value = 6000

# --- Jason Schema: same Rate (F) / Total (G) population for its rows, and
#     the corrected employee id (column B) ---
for ($r = 2; $r -le 6; $r++) {
    $jsSheet.Cells.Item($r, 2).Value = "emp_35u1tnme"
    $jsSheet.Cells.Item($r, 6).Value = 150
    $jsSheet.Cells.Item($r, 7).Value = 1200
}
